# RMI files through Dec 1
# Update GDP Growth Rates workbook to reflect the November STEO data
# (replacing the July STEO / May 2020 pandemic-impact figures).

$wb = $excel.ActiveWorkbook

# Turn off iterative calculation (the workbook no longer needs it) -
# equivalent to clearing Formulas > Calculation Options > Enable Iterative
# Calculation in the UI.
$excel.Iteration = $false

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Data")

# --- Data sheet: source label + updated GDP figures (update text first so new
# shared strings are interned in the same order Excel would have written them) ---
$wsData.Range("A3").Value = "November STEO"
$wsData.Range("B3").Value = 19092
$wsData.Range("C3").Value = 18411
$wsData.Range("D3").Value = 19098

# --- About sheet text updates ---
$wsAbout.Range("B6").Value = "January 2020 and November 2020"
$wsAbout.Range("A27").Value = "As of EPS 3.1, this variable is set up to model the impacts of the 2020"
$wsAbout.Range("A28").Value = "SARS-CoV-2 pandemic.  It uses the latest data available as of November 10,"

# Update the selected cell on the Data sheet
$wsData.Range("B12").Select()

# Update the selected cell on the About sheet (leave it the active sheet/selection)
$wsAbout.Range("A29").Select()

$wb.Save()
